# Updated cryptos list on Wed Jul 19 09:06:53 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table
# with newly scraped values. Rows 47/48 additionally swap rank between
# "EnergySwap" and "Aptos" (Coin/Link/Price/Volume all change).
#
# Price-looking strings (e.g. "1.002", "241.52") are pre-formatted as
# Text ("@") before assignment so Excel keeps them as literal strings
# instead of silently converting them to numbers - matching how the
# source data already stores them as text in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.959.47"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.907.45"
$ws.Range("E3").Value = "  +0.23%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5 - XRP
Set-TextValue "D5" "0.7763"
$ws.Range("E5").Value = "  +4.19%  "

# Row 6 - BNB
Set-TextValue "D6" "241.52"
$ws.Range("E6").Value = "  -0.30%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3139"
$ws.Range("E8").Value = "  +2.19%  "

# Row 9 - Solana
Set-TextValue "D9" "25.87"
$ws.Range("E9").Value = "  +0.90%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.06854"
$ws.Range("E10").Value = "  -0.84%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07962"
$ws.Range("E11").Value = "  -1.07%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.903.16"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13 - Polygon
Set-TextValue "D13" "0.7372"
$ws.Range("E13").Value = "  -2.66%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.178"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15 - Litecoin
Set-TextValue "D15" "92.42"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "29.972.08"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17 - Avalanche
Set-TextValue "D17" "13.84"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18 - Uniswap
Set-TextValue "D18" "5.840"
$ws.Range("E18").Value = "  -5.55%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "244.26"
$ws.Range("E19").Value = "  +2.76%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.000007690"
$ws.Range("E20").Value = "  -1.22%  "

# Row 21 - Dai (price only, Volume unchanged)
Set-TextValue "D21" "1.002"

# Row 22 - WrappedliquidstakedEther2.0
Set-TextValue "D22" "2.147.54"
$ws.Range("E22").Value = "  -0.35%  "

# Row 23 - BinanceUSD
Set-TextValue "D23" "1.003"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24 - Chainlink
Set-TextValue "D24" "6.829"
$ws.Range("E24").Value = "  -3.64%  "

# Row 25 - Monero
Set-TextValue "D25" "168.63"
$ws.Range("E25").Value = "  +0.48%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.246"
$ws.Range("E26").Value = "  -0.81%  "

# Row 27 - Stellar
Set-TextValue "D27" "0.1361"
$ws.Range("E27").Value = "  +7.37%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "18.84"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29 - LidoDAOToken
Set-TextValue "D29" "2.015"
$ws.Range("E29").Value = "  -1.81%  "

# Row 30 - Toncoin
Set-TextValue "D30" "1.369"
$ws.Range("E30").Value = "  +1.30%  "

# Row 31 - PancakeSwap (Volume only, Price unchanged)
$ws.Range("E31").Value = "  -0.73%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.302"
$ws.Range("E32").Value = "  -0.03%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "4.062"
$ws.Range("E33").Value = "  +0.32%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.05454"
$ws.Range("E34").Value = "  +2.87%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.247"
$ws.Range("E35").Value = "  -3.24%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.7281"
$ws.Range("E36").Value = "  -1.72%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.735"
$ws.Range("E37").Value = "  +0.33%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.01924"
$ws.Range("E38").Value = "  -1.09%  "

# Row 39 - MXToken
Set-TextValue "D39" "2.785"
$ws.Range("E39").Value = "  +0.77%  "

# Row 40 - FraxShare
Set-TextValue "D40" "6.107"
$ws.Range("E40").Value = "  -2.41%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.4393"
$ws.Range("E41").Value = "  -1.57%  "

# Row 42 - Aave
Set-TextValue "D42" "71.42"
$ws.Range("E42").Value = "  -1.80%  "

# Row 43 - PaxDollar
Set-TextValue "D43" "1.003"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44 - TrustWalletToken
Set-TextValue "D44" "0.8399"
$ws.Range("E44").Value = "  +1.00%  "

# Row 45 - RenderToken
Set-TextValue "D45" "1.860"
$ws.Range("E45").Value = "  -4.90%  "

# Row 46 - Quant
Set-TextValue "D46" "99.95"
$ws.Range("E46").Value = "  -1.36%  "

# Rows 47/48 swap rank between EnergySwap and Aptos, with refreshed data
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.467"
$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.704"
$ws.Range("E48").Value = "  -0.50%  "

# Row 49 - Maker
Set-TextValue "D49" "975.30"
$ws.Range("E49").Value = "  +7.95%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "2.056.22"
$ws.Range("E50").Value = "  -0.04%  "

# Row 51 - Elrond
Set-TextValue "D51" "36.03"
$ws.Range("E51").Value = "  -1.58%  "
